$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.930.80'
$ws.Range('E2').Value = '  +1.02%  '
$ws.Range('D3').Value = '2.237.58'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.62'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.28%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.570'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.51%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.537'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.43'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0821'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.38'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.55%  '
$ws.Range('E13').Value = '  -2.58%  '
$ws.Range('D14').Value = '2.579.30'
$ws.Range('E14').Value = '  -0.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.842'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.10%  '
$ws.Range('D16').Value = '2.243.20'
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.04'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.16%  '
$ws.Range('D18').Value = '43.840.45'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.92'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -9.34%  '
$ws.Range('D20').Value = '0.0₃0965'
$ws.Range('E20').Value = '  -1.81%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.36'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '64.94'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '234.03'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.36%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.27'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.84'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.05'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '158.75'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.96'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.85%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0836'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.49%  '
$ws.Range('E34').Value = '  -0.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.15'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.110'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.88'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.117'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.86'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.62'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.02'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0310'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.47%  '
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').Value = '1.729.82'
$ws.Range('E44').Value = '  -5.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.194'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '80.73'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '73.33'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.11'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '101.51'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.55%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.64'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.36%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '56.87'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.65%  '
